$d = $word.ActiveDocument

# --- Change 1: remove the (now-extraneous) empty paragraph that followed the
# "probe implantation ... lesioning/histology." paragraph. There were two
# consecutive empty paragraphs there; delete the first of the pair so only
# one remains before "To Incorporate Intan Electrophysiological Data ...".
$found = $d.Content.Find.Execute("probe implantation (to verify probe works)")
if (-not $found) {
    throw "Could not find the 'probe implantation' paragraph"
}
$hostPara = $d.Content.Paragraphs.First
# Find.Execute leaves $d.Content collapsed-ish; re-resolve via a fresh search range
$searchRange = $d.Content
$searchRange.Find.Execute("probe implantation (to verify probe works)") | Out-Null
$ownerPara = $searchRange.Paragraphs(1)
$emptyPara = $ownerPara.Next()
if ($emptyPara.Range.Text.Trim().Length -ne 0) {
    throw "Expected an empty paragraph after the 'probe implantation' paragraph"
}
$emptyPara.Range.Delete()

# --- Change 2: drop the <w:lastRenderedPageBreak/> marker that currently sits
# on the "Must have the path in the 'current folder' ..." run. Re-submitting
# the paragraph's own WordOpenXML through InsertXML rewrites the paragraph
# without that (non-round-tripping) marker while preserving every other
# attribute untouched.
$rng1 = $d.Content
$rng1.Find.Execute("Must have the path in the") | Out-Null
$p1 = $rng1.Paragraphs(1)
$p1.Range.InsertXML($p1.Range.WordOpenXML)

# --- Change 3: add <w:lastRenderedPageBreak/> onto the run that starts the
# "nexData = intan2nex..." paragraph (i.e. move the page-break marker down
# onto this paragraph instead).
$rng2 = $d.Content
$rng2.Find.Execute("nexData = intan2nex") | Out-Null
$p2 = $rng2.Paragraphs(1)
$marker = "<w:t>nexData = intan2nex</w:t>"
$xml2 = $p2.Range.WordOpenXML
if ($xml2.IndexOf($marker) -lt 0) {
    throw "Could not locate the 'nexData = intan2nex' run to tag"
}
$newXml2 = $xml2.Replace($marker, "<w:lastRenderedPageBreak/>" + $marker)
$p2.Range.InsertXML($newXml2)
